$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Año (Year) and Sem (Week) values for rows 2-5
$ws.Range("B2").Value = 2025
$ws.Range("C2").Value = 30

$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 30

$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 30

$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 30

# Update execution time in row 8
$ws.Range("J8").Value = "17:00:16"
